$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 272) holds a date serial number that was
# bumped by one day (46060 -> 46061) in the latest export.
$ws.Range("C2:C272").Value = 46061
